# Fix a student's first-name entry that was caught during inspection, then
# sort the paper list by publication Date (ascending) as originally intended.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Correct "Katrina Lee" -> "Katrina" (surname had been mistakenly typed into the first-name cell)
$ws.Range("B2").Value = "Katrina"

# 2. Sort the data range (A2:H12) ascending by the Date column (G)
$dataRange = $ws.Range("A2:H12")
$sortKey = $ws.Range("G2:G12")
$dataRange.Sort($sortKey, 1)

# Record the sort on the worksheet's Sort object so the sort state is persisted
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("G2:G12"))
$sortObj.SetRange($ws.Range("A2:H12"))
$sortObj.Header = 2
$sortObj.Apply()

# 3. Widen the two newly-visible columns (G, H) to match the rest of the table
$ws.Columns.Item(7).ColumnWidth = 11.498697916666666
$ws.Columns.Item(8).ColumnWidth = 18.666666666666668

# 4. Leave the selection on the sorted block
$ws.Range("B2:H12").Select()
